$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 02df0fce... row
# refreshed to a newer handback-generation timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-13 10:59:35"

# zh-cn handback table: 02df0fce... row got a fresh handoff/handback cycle.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-13 10:59:27"
$wsZhCn.Range("K2").Value = "2016-08-13 10:59:55"

# de-de handback table: 02df0fce... row got a fresh handoff/handback cycle.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-13 10:59:35"
$wsDeDe.Range("K2").Value = "2016-08-13 11:00:12"
